# Update term package metadata to version 1.1.0, per commit "Added 1.1.0 of term".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds property names, column B holds the corresponding values.
# Row 3 = Version, Row 8 = Date (see sharedStrings ordering: Property/Value,
# URL/..., Version/1.0.0, Name/..., Title/..., Status/..., Experimental/...,
# Date/2023-06-07T11:52:14+02:00, ...).
$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
